$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header cell's formatting (bold, centered, bordered) to the new header cell
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("H1").Value = "Save"

# New "Save" column values
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 1
$ws.Range("H6").Value = 0
